# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 50 (pushing existing rows 50-76 down to 51-77)
# and populate the newly inserted row with the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 44460
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = 100112031
$ws.Cells.Item(50, 7).Value = "Poroto verde"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 20
$ws.Cells.Item(50, 11).Value = 40000
$ws.Cells.Item(50, 12).Value = 40000
$ws.Cells.Item(50, 13).Value = 40000
$ws.Cells.Item(50, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(50, 15).Value = "Perú"
$ws.Cells.Item(50, 16).Value = 1600
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
